$d = $word.ActiveDocument

# The diff appends two brand-new paragraphs right after the paragraph that
# contains "是小钟雷的生日" (the last paragraph in the document), and the
# trailing _GoBack bookmark (bookmarkStart/bookmarkEnd) moves along with the
# edit point, ending up after the very last new paragraph.
$anchorText = "是小钟雷的生日"
$newText1   = "今天是2022年9月17日"
$newText2   = "星期六，又在上课了"

# Use Find/Replace to append the new text onto the end of the anchor
# paragraph (as one run). Editing the content this way - rather than
# inserting a brand-new empty paragraph at the very end of the document -
# keeps the document's "last edit" location (and therefore the _GoBack
# bookmark) tracking with this edit, so it ends up in the right place once
# the paragraph gets split below.
$combined = $anchorText + $newText1 + $newText2
$found = $d.Content.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, $combined, 2)

# The paragraph that used to hold just $anchorText now holds $combined, and
# (since the document has no content after it) is still the last paragraph.
$p = $d.Paragraphs($d.Paragraphs.Count)
$pStart = $p.Range.Start

# Split "...是小钟雷的生日|今天是2022年9月17日星期六，又在上课了" into two
# paragraphs by inserting a paragraph mark right after the anchor text.
$split1 = $pStart + $anchorText.Length
$r1 = $d.Range($split1, $split1)
$r1.InsertParagraphAfter()

# Split "...今天是2022年9月17日|星期六，又在上课了" into two paragraphs the
# same way (offset shifts by 1 for the paragraph mark just inserted).
$split2 = $split1 + 1 + $newText1.Length
$r2 = $d.Range($split2, $split2)
$r2.InsertParagraphAfter()
